$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# NOTE: Range.InsertXML replaces the whole paragraph(s) spanned by the range
# it is called on (this engine resolves it at paragraph granularity, not a
# literal character splice). So whenever we "insert" new paragraphs next to
# an existing one, we re-supply that existing paragraph's own content as
# part of the replacement XML together with the new paragraph(s), using the
# existing paragraph as the anchor/host for the call.

# --- 1) Split "figma" into its own run with spell-check proofErr markers ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*translating figma design files*") {
        $xml1 = '<w:p ' + $wNs + '>' +
                  '<w:pPr><w:pStyle w:val="Undertittel"/></w:pPr>' +
                  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
                  '<w:r><w:br/></w:r>' +
                  '<w:r><w:t xml:space="preserve">At the start of this assignment I felt confident about the design but aware that translating </w:t></w:r>' +
                  '<w:proofErr w:type="spellStart"/>' +
                  '<w:r><w:t>figma</w:t></w:r>' +
                  '<w:proofErr w:type="spellEnd"/>' +
                  '<w:r><w:t xml:space="preserve"> design files into structured, semantic HTML would be a challenge. The design includes multiple page types such as homepage, product listing, product detail pages and a checkout flow, which requires careful planning to keep the structure consistent and accessible.</w:t></w:r>' +
                '</w:p>'
        $p.Range.InsertXML($xml1)
        break
    }
}

# --- 2) Split "lated" into its own run with spell-check proofErr markers ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*styling lated, instead*") {
        $xml2 = '<w:p ' + $wNs + '>' +
                  '<w:r><w:t xml:space="preserve">I expected the workload to be manageable but time consuming, especially when ensuring that the website should have semantic HTML, accessibility and have good responsiveness across all pages. To stay organized I am planning to focus on structure first, and styling </w:t></w:r>' +
                  '<w:proofErr w:type="spellStart"/>' +
                  '<w:r><w:t>lated</w:t></w:r>' +
                  '<w:proofErr w:type="spellEnd"/>' +
                  '<w:r><w:t>, instead of trying to perfect everything at once.</w:t></w:r>' +
                '</w:p>'
        $p.Range.InsertXML($xml2)
        break
    }
}

# --- 3) Insert new reflection paragraphs about the cart badge, right before
#        "Inspiration & Sources" (anchored on the preceding paragraph so its
#        own content is preserved). ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Inspiration & Sources*") {
        $prev = $p.Previous()
        $xml3 = '<w:p ' + $wNs + '><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
                  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>I spent some time comparing my implementation directly to the Figma design I created earlier. I adjusted spacing alignment and visual hierarchy to achieve closer match instead of restructuring the markup.</w:t></w:r>' +
                '</w:p>' +
                '<w:p ' + $wNs + '><w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p>' +
                '<w:p ' + $wNs + '>' +
                  '<w:pPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
                  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">During the process of making the </w:t></w:r>' +
                  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>checkout I wanted to make a badge for the shopping cart. With inspiration from komplett.no and help from ChatGPT I managed to create something I am very happy with.</w:t></w:r>' +
                '</w:p>'
        $prev.Range.InsertXML($xml3)
        break
    }
}

# --- 4) Insert the two new source paragraphs right after "Inspiration &
#        Sources" (anchored on that paragraph itself, re-supplying its own
#        heading content). ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Inspiration & Sources*") {
        $xml4 = '<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="Overskrift1"/></w:pPr><w:r><w:t>Inspiration &amp; Sources</w:t></w:r></w:p>' +
                '<w:p ' + $wNs + '><w:r><w:t>Komplett.no</w:t></w:r></w:p>' +
                '<w:p ' + $wNs + '><w:r><w:t>ChatGPT by OpenAI</w:t></w:r></w:p>'
        $p.Range.InsertXML($xml4)
        break
    }
}
